# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------
# 1) Swap the display order of a few country name / data pairs.
#    The rows keep their row number, but the country label (and the data
#    that "belongs" to that country) move to the other row of the pair.
# ---------------------------------------------------------------------------

# Curazao (row 198) <-> Fiyi (row 199)
$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# Montserrat (row 210) <-> Seychelles (row 211)
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Islas Virgenes Britanicas (row 213) <-> Papua Nueva Guinea (row 214)
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# ---------------------------------------------------------------------------
# 2) Update the "last updated" timestamp banner.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 11:05"

# ---------------------------------------------------------------------------
# 3) Refresh the statistics for a handful of countries / regions.
# ---------------------------------------------------------------------------

# Row 25 - Banglades
$ws.Range("B25").Value = 44608
$ws.Range("C25").Value = 1764
$ws.Range("D25").Value = 9375
$ws.Range("E25").Value = 34623
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 610

# Row 41 - Rumania
$ws.Range("E41").Value = 4900
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 1253

# Row 101 - Sri Lanka
$ws.Range("D101").Value = 781
$ws.Range("E101").Value = 767

# Row 112 - Hong Kong
$ws.Range("B112").Value = 1083
$ws.Range("C112").Value = 3
$ws.Range("D112").Value = 1036
$ws.Range("E112").Value = 43

# Row 141 - Estado de Palestina
$ws.Range("B141").Value = 447
$ws.Range("C141").Value = 1
$ws.Range("E141").Value = 76
